# Projects-Information.xlsx — update the Electrical row's image path from the
# old Flask-style `url_for(...)` reference to a plain `/images/...` path, and
# leave the view focused on that cell (scrolled so column C is left-most,
# with E2 selected), matching how the author left the sheet after editing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# E2 ("Image path" column, Electrical row): swap the old templated path for
# the new static path.
$ws.Range("E2").Value = "/images/it in the hall.jpg"

# Leave the sheet scrolled/selected the way it was when saved: view starting
# at column C, with E2 as the active cell.
$excel.ActiveWindow.ScrollColumn = 3
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E2").Select()
